$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja1: update existing rows 2-3 and fill rows 4-6 with new OT data ---
$ws1.Range("C2").Value = "MPI_ANUAL_01-20_UNILEVER"
$ws1.Range("D2").Value = "REDCOM01S26550"

$ws1.Range("C3").Value = "MPI_ANUAL_01-20_ESTRUMENTAL"
$ws1.Range("D3").Value = "REDCOM01S26551"

$ws1.Range("C4").Value = "MPI_ANUAL_01-20_JARAMILLO MORA"
$ws1.Range("D4").Value = "MLL_TCAL_540681_024"

$ws1.Range("C5").Value = "MPI_ANUAL_01-20_PRONAVICOLA MEDIACANOA"
$ws1.Range("D5").Value = "UL9907_METROCELDA_PRONAVICOLA_MEDIOCANOA_7"

$ws1.Range("C6").Value = "MPI_SEMESTRAL_01-20_FLORA INDUSTRIAL (TRASLADO)"
$ws1.Range("D6").Value = "UL088"

# --- Hoja1: extend sheet down to row 62 with the same blank formatted rows ---
$ws1.Range("A21:E21").Copy() | Out-Null
$ws1.Range("A22:E62").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Hoja1: update the selected cell ---
$ws1.Range("C10").Select() | Out-Null

# --- Hoja2: update existing rows 2-3 and add new rows 4-6 ---
$ws2.Range("A2").Value = "S26550"
$ws2.Range("B2").Value = "REDCOM01S26550"
$ws2.Range("C2").Value = 5182432
$ws2.Range("E2").Value = "MPI_ANUAL_01-20_UNILEVER"

$ws2.Range("A3").Value = "S26551"
$ws2.Range("B3").Value = "REDCOM01S26551"
$ws2.Range("C3").Value = 5182433
$ws2.Range("E3").Value = "MPI_ANUAL_01-20_ESTRUMENTAL"

$ws2.Range("A4").Value = "S25809"
$ws2.Range("B4").Value = "MLL_TCAL_540681_024"
$ws2.Range("C4").Value = 5182434
$ws2.Range("E4").Value = "MPI_ANUAL_01-20_JARAMILLO MORA"

$ws2.Range("A5").Value = "S26909"
$ws2.Range("B5").Value = "UL9907_METROCELDA_PRONAVICOLA_MEDIOCANOA_7"
$ws2.Range("C5").Value = 5182435
$ws2.Range("E5").Value = "MPI_ANUAL_01-20_PRONAVICOLA MEDIACANOA"

$ws2.Range("A6").Value = "S27059"
$ws2.Range("B6").Value = "UL088"
$ws2.Range("C6").Value = 5182436
$ws2.Range("E6").Value = "MPI_SEMESTRAL_01-20_FLORA INDUSTRIAL (TRASLADO)"

# --- Hoja2: update the selected cell ---
$ws2.Range("E14").Select() | Out-Null

# --- Make Hoja2 the active sheet (matches activeTab="1" / tabSelected) ---
$ws2.Activate() | Out-Null

# --- Restore window size/position to match the saved workbook view ---
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 10420
